# Saldo.xlsx — "Export" sheet maintenance edit
#
# Target change (per the author's diff):
#   - Remove the row for account 004572740 / PAULO / 231567.92
#     (was directly under the header, Excel row 2).
#   - Insert a new row for account 005152037 / RODRIGO / 25057.12
#     immediately above the THIAGO (005064129) row.
#   - Change DANIELA's (004001621) Saldo from 20000 to 10917.49.
#   - Remove the row for account 004452476 / IVONE / 19937.62
#     (was directly under DANIELA's row).
#
# All other rows stay exactly as-is and simply shift position as a
# consequence of the row delete/insert operations above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Delete the PAULO row (account 004572740), originally row 2 ---------
$ws.Rows(2).Delete()

# After this delete, the rows directly under the header are now:
#   row 2 = FERNANDO (005726697)
#   row 3 = THIAGO   (005064129)
#   row 4 = DANIELA  (004001621)
#   row 5 = IVONE    (004452476)

# --- 2) Insert the new RODRIGO row right above THIAGO (row 3) --------------
$ws.Rows(3).Insert()

# Account numbers are text with significant leading zeros, so force the
# cell to text (leading apostrophe) and then strip the resulting
# "number stored as text" formatting so the cell carries no explicit
# style, matching the rest of the sheet's plain data cells.
$ws.Range("A3").Value = "'005152037"
$ws.Range("A3").ClearFormats()
$ws.Range("B3").Value = "RODRIGO"
$ws.Range("C3").Value = 25057.12

# After the insert, the rows are now:
#   row 2 = FERNANDO
#   row 3 = RODRIGO (new)
#   row 4 = THIAGO
#   row 5 = DANIELA
#   row 6 = IVONE

# --- 3) Update DANIELA's Saldo value ---------------------------------------
$ws.Range("C5").Value = 10917.49

# --- 4) Delete the IVONE row -------------------------------------------------
$ws.Rows(6).Delete()
